$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new card ("Mirari's Wake") was bought/tracked but the price-fetch
# commit never got pulled in, so it's missing from the album. Insert a
# new row for it right after the current row 58 (before the old row 59),
# pushing every row below it down by one - matches how a new
# entry gets added to this running list.
$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value = "Mirari's Wake"
$ws.Range("B59").Value = "Modern Horizons 2"
$ws.Range("C59").Value = "V.2"
$ws.Range("D59").Value = 2.79

# Match the author's on-screen scroll/selection state after the edit.
$null = $ws.Range("A79").Select()
